# Apply the edit described by the commit:
# "Fixed bugs and added print statements to console output"
#
# Net effect on the workbook data (the author re-ran the cost model for a
# different shop/site - BME_BCCW instead of BME_DI_VGH - which changed the
# title, the overhead/rate inputs, and the asset list, and one of the two
# previous asset rows was dropped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet / report title ------------------------------
$ws.Name = "BME_BCCW"
$ws.Range("A1").Value = "BME_BCCW: Annual Service Delivery Costs for Net New Equipment"

# --- Updated OH Information figures ------------------------------------
$ws.Range("B4").Value = 1336293.590552836
$ws.Range("B5").Value = 561058.118
$ws.Range("B6").Value = 690897.0139690499
$ws.Range("B7").Value = 84338.45858378569

# --- Updated Rates ------------------------------------------------------
$ws.Range("B10").Value = 51.30578674124474
$ws.Range("B11").Value = 39.80729857819905

# --- Updated (single remaining) asset row -------------------------------
$ws.Range("A15").Value = "PHSA"
$ws.Range("B15").Value = "BCCW"
$ws.Range("C15").Value = "BCCH"
$ws.Range("D15").Value = 12378
$ws.Range("E15").Value = "ELECTRIC SIGNAL AMPLIFIERS"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.046341463414635

# --- Drop the second asset row that no longer applies -------------------
$ws.Rows.Item(16).Delete()
